$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# This shared string is used by:
#   Overview sheet: E2 (zh-cn status), F2 (de-de status)
#   zh-cn sheet:    C2 (Status)
#   de-de sheet:    C2 (Status)
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: put the string literal on the left of -eq so a Boolean cell
        # value (e.g. True/False) is not coerced into matching a non-empty
        # string (PowerShell comparison-operand-type coercion gotcha).
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# --- Narrow the "Status" columns ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

# zh-cn sheet: column C (Status)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511

# de-de sheet: column C (Status)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
